# "two more permutations for adding to cart"
# Adds a new permutation row to the Smartwatches sheet and makes it the
# active sheet/selection (mirrors what Excel records when a user enters
# data on that sheet and saves).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Smartwatches")
$ws.Activate()

$ws.Range("A4").Value = "smartwatches"
$ws.Range("B4").Value = "6"
$ws.Range("C4").Value = "0"
$ws.Range("D4").Value = "0"
$ws.Range("E4").Value = "11102"
$ws.Range("F4").Value = "has been added to cart."

# Match formatting of the row above (Menlo font style) instead of the
# default numeric/text style that a plain Value assignment would pick up.
$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)

$ws.Range("F4").Select()
